$d = $word.ActiveDocument

# 1. Update building ID in title (avoid touching the apostrophe so Word's
#    smart-quote autocorrect doesn't turn it into a curly quote)
$d.Content.Find.Execute("Immeuble ID: 1", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Immeuble ID: 4", 2)

# 2. Update owner name
$d.Content.Find.Execute("Luigi Bros", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Luigi Brothers", 2)

# 3. Update contact info (remove email, change phone number)
$d.Content.Find.Execute("Contact : luigi@nintendo.com, 07 47 58 69 47", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Contact : , 07 58 47 61 25", 2)

# 4. Update total rental income
$d.Content.Find.Execute("Revenus locatifs totaux : 750.0 €", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Revenus locatifs totaux : 8300.0 €", 2)

# 5. Update total expenses
$d.Content.Find.Execute("Dépenses totales (charges et travaux) : 253.4 €", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Dépenses totales (charges et travaux) : 108.4 €", 2)

# 6. Update unpaid taxes/bills
$d.Content.Find.Execute("Taxes et factures impayées : 0.0 €", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Taxes et factures impayées : 2400.0 €", 2)
